$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDA")

# Flatten the named range _b (EDA!$H$3:$H$20) into a row, spilling M3:AD3
$ws.Range("M3:AD3").FormulaArray = "=TOROW(_b)"

# Copy _b down as a column, spilling L4:L21
$ws.Range("L4:L21").FormulaArray = "=_b"

# Build the "same value" comparison matrix, spilling M4:AD21
$ws.Range("M4:AD21").FormulaArray = "=N(_b=TOROW(_b))"

$ws.Range("L5").Select()
